# Misc updates - mostly so MASTER excel can be read in directly from UCSRB server
#
# 1) Two header columns renamed (M1, N1)
# 2) Entiat River Lake 03 (row 6): Riparian-Disturbance_score / Riparian_Mean / HQ_Sum / HQ_Pct recalculated
# 3) Entiat River Lake 05 (row 8): Bull.Trout.Reach flag flipped no -> yes
# 4) A new reach "Methow River Fawn 02" inserted as row 17, pushing the remaining
#    Methow/Nason/Twisp/White reaches down one row, with several of those rows'
#    scores also refreshed from the updated MASTER source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param(
        [int]$RowNum,
        [object[]]$Values
    )
    $arr = New-Object 'object[,]' 1,$Values.Length
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $arr[0,$i] = $Values[$i]
    }
    $endCol = [char]([int][char]'A' + $Values.Length - 1)
    $rangeAddr = "A$RowNum`:$endCol$RowNum"
    $ws.Range($rangeAddr).Value = $arr
}

# --- Header renames ---
$ws.Range("M1").Value = "FloodplainConnectivity_score"
$ws.Range("N1").Value = "Off-Channel/Side-Channels_score"

# --- Row 6: Entiat River Lake 03 score refresh ---
$ws.Range("Q6").Value = 5
$ws.Range("R6").Value = 4
$ws.Range("T6").Value = 40
$ws.Range("U6").Value = 0.8888888888888888

# --- Row 8: Entiat River Lake 05 Bull Trout reach flag ---
$ws.Range("F8").Value = "yes"

# --- Rows 17-30: insert "Methow River Fawn 02" and shift/refresh the rest ---
Set-RowValues 17 @("Methow River Fawn 02", "Methow", "Methow River-Fawn Creek", "yes", "yes", "yes", 3, 5, 4, 5, 3, 5, 3, 3, 5, 1, 5, 3, 1, 32, 0.7111111111111111, 5, 3)
Set-RowValues 18 @("Methow River Fawn 04", "Methow", "Methow River-Fawn Creek", "yes", "yes", "yes", 5, 5, 5, 5, 3, 3, 5, 5, 5, 3, 5, 4, 1, 36, 0.8, 3, 3)
Set-RowValues 19 @("Methow River Rattlesnake 01", "Methow", "Methow River-Rattlesnake Creek", "yes", "yes", "yes", 3, 3, 3, 5, 3, 1, 3, 5, 5, 1, 3, 2, 5, 32, 0.7111111111111111, 5, 3)
Set-RowValues 20 @("Methow River Rattlesnake 05", "Methow", "Methow River-Rattlesnake Creek", "yes", "yes", "yes", 5, 5, 5, 5, 1, 5, 5, 5, 3, 3, 3, 3, 5, 37, 0.8222222222222222, 3, 3)
Set-RowValues 21 @("Methow River Rattlesnake 06", "Methow", "Methow River-Rattlesnake Creek", "yes", "yes", "yes", 5, 5, 5, 5, 5, 5, 5, 5, 5, 5, 5, 5, 5, 45, 1, 1, 5)
Set-RowValues 22 @("Methow River Thompson 07", "Methow", "Methow River-Thompson Creek", "yes", "yes", "yes", 3, 5, 4, 5, 5, 1, 3, 3, 5, 1, 3, 2, 5, 33, 0.7333333333333333, 5, 3)
Set-RowValues 23 @("Methow River Thompson 08", "Methow", "Methow River-Thompson Creek", "yes", "yes", "yes", 3, 5, 4, 5, 5, 1, 3, 3, 5, 1, 3, 2, 5, 33, 0.7333333333333333, 5, 3)
Set-RowValues 24 @("Nason Creek Lower 01", "Wenatchee", "Lower Nason Creek", "yes", "yes", "yes", 3, 5, 4, 5, 5, 5, 3, 5, 5, 3, 5, 4, 1, 37, 0.8222222222222222, 3, 3)
Set-RowValues 25 @("Nason Creek Lower 02", "Wenatchee", "Lower Nason Creek", "yes", "yes", "yes", 3, 5, 4, 5, 5, 5, 3, 5, 5, 5, 5, 5, 3, 40, 0.8888888888888888, 3, 3)
Set-RowValues 26 @("Nason Creek Lower 03", "Wenatchee", "Lower Nason Creek", "yes", "yes", "yes", 3, 5, 4, 5, 5, 3, 3, 5, 5, 3, 5, 4, 1, 35, 0.7777777777777778, 5, 3)
Set-RowValues 27 @("Twisp River Middle 01", "Methow", "Middle Twisp River", "yes", "yes", "yes", 3, 3, 3, 5, 5, 3, 3, 5, 5, 3, 3, 3, 1, 33, 0.7333333333333333, 5, 3)
Set-RowValues 28 @("Twisp River Middle 02", "Methow", "Middle Twisp River", "yes", "yes", "yes", 3, 3, 3, 5, 5, 5, 3, 5, 5, 3, 3, 3, 1, 35, 0.7777777777777778, 5, 3)
Set-RowValues 29 @("Twisp River Middle 06", "Methow", "Middle Twisp River", "yes", "yes", "yes", 5, 5, 5, 5, 1, 5, 5, 5, 5, 1, 5, 3, 3, 37, 0.8222222222222222, 3, 3)
Set-RowValues 30 @("White River Lower 08", "Wenatchee", "Lower White River", "yes", "yes", "yes", 5, 5, 5, 5, 3, 5, 5, 1, 1, 3, 5, 4, 5, 34, 0.7555555555555555, 5, 3)
